$d = $word.ActiveDocument
$sec = $d.Sections.First

# Headers: Item(1) -> header2.xml (default, id=3), Item(2) -> header1.xml (first page, id=1)
# Footers: Item(1) -> footer2.xml (default, id=4), Item(2) -> footer1.xml (first page, id=2)

$hDefault = $sec.Headers.Item(1)
[void]$hDefault.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

$hFirst = $sec.Headers.Item(2)
[void]$hFirst.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

$fDefault = $sec.Footers.Item(1)
[void]$fDefault.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

$fFirst = $sec.Footers.Item(2)
[void]$fFirst.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

Write-Host "Renamed inline shapes in headers/footers."
